$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.014515734379126
$ws.Range("D2").Value = 1.017160917744017
$ws.Range("E2").Value = 1.016247304167528
$ws.Range("F2").Value = 1.025669361290694
$ws.Range("I2").Value = 1.025718333119251
$ws.Range("J2").Value = 1.019746526450355
$ws.Range("K2").Value = 1.020010933592147
$ws.Range("L2").Value = 1.019100044737615
$ws.Range("M2").Value = 1.028494255553923
$ws.Range("N2").Value = 1.021194683484476
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.015942431034938
$ws.Range("D3").Value = 1.018499098826531
$ws.Range("E3").Value = 1.017473029416238
$ws.Range("F3").Value = 1.027163534987985
$ws.Range("I3").Value = 1.025712513985777
$ws.Range("J3").Value = 1.020805312090053
$ws.Range("K3").Value = 1.021153178957294
$ws.Range("L3").Value = 1.020129944729006
$ws.Range("M3").Value = 1.029793918524472
$ws.Range("N3").Value = 1.022254972721227
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.016864741584167
$ws.Range("D4").Value = 1.019364497503446
$ws.Range("E4").Value = 1.018265668587749
$ws.Range("F4").Value = 1.028127783493812
$ws.Range("I4").Value = 1.02570617741944
$ws.Range("J4").Value = 1.021489217222405
$ws.Range("K4").Value = 1.021891272667994
$ws.Range("L4").Value = 1.020795324697716
$ws.Range("M4").Value = 1.030631820601973
$ws.Range("N4").Value = 1.022939849077312
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.017252281974536
$ws.Range("D5").Value = 1.0197281979913
$ws.Range("E5").Value = 1.018598781937172
$ws.Range("F5").Value = 1.028532541319648
$ws.Range("I5").Value = 1.025702897401982
$ws.Range("J5").Value = 1.021776448032499
$ws.Range("K5").Value = 1.022201328780118
$ws.Range("L5").Value = 1.021074806907096
$ws.Range("M5").Value = 1.030983345366157
$ws.Range("N5").Value = 1.023227487788102
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.017317340227731
$ws.Range("D6").Value = 1.019789258419248
$ws.Range("E6").Value = 1.018654706685745
$ws.Range("F6").Value = 1.028600466205797
$ws.Range("I6").Value = 1.02570231052569
$ws.Range("J6").Value = 1.021824658886026
$ws.Range("K6").Value = 1.022253374714305
$ws.Range("L6").Value = 1.021121719036624
$ws.Range("M6").Value = 1.031042325298096
$ws.Range("N6").Value = 1.02327576710657
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.016869920691797
$ws.Range("D7").Value = 1.0193693577253
$ws.Range("E7").Value = 1.018270120099869
$ws.Range("F7").Value = 1.028133194286379
$ws.Range("I7").Value = 1.025706136013405
$ws.Range("J7").Value = 1.021493056320814
$ws.Range("K7").Value = 1.02189541658435
$ws.Range("L7").Value = 1.020799060104593
$ws.Range("M7").Value = 1.030636520552732
$ws.Range("N7").Value = 1.022943693627681
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.014998074589893
$ws.Range("D8").Value = 1.017613267686982
$ws.Range("E8").Value = 1.016661647274152
$ws.Range("F8").Value = 1.026174862394676
$ws.Range("I8").Value = 1.025716898558997
$ws.Range("J8").Value = 1.020104599114523
$ws.Range("K8").Value = 1.020397173793071
$ws.Range("L8").Value = 1.019448320848024
$ws.Range("M8").Value = 1.028934120712789
$ws.Range("N8").Value = 1.021553264652899
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.011692789586293
$ws.Range("D9").Value = 1.014514798357592
$ws.Range("E9").Value = 1.013823382982949
$ws.Range("F9").Value = 1.022703994053582
$ws.Range("I9").Value = 1.025716195644615
$ws.Range("J9").Value = 1.017648573175214
$ws.Range("K9").Value = 1.017749094304407
$ws.Range("L9").Value = 1.017060042384645
$ws.Range("M9").Value = 1.025910564019418
$ws.Range("N9").Value = 1.019093750875045
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.009484258547581
$ws.Range("D10").Value = 1.012446127802575
$ws.Range("E10").Value = 1.011928283348431
$ws.Range("F10").Value = 1.020376235852606
$ws.Range("I10").Value = 1.02570253460563
$ws.Range("J10").Value = 1.016004634410336
$ws.Range("K10").Value = 1.015978047787012
$ws.Range("L10").Value = 1.015462151133398
$ws.Range("M10").Value = 1.023878601063714
$ws.Range("N10").Value = 1.017447477528559
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.008526662915644
$ws.Range("D11").Value = 1.011549580244184
$ws.Range("E11").Value = 1.011106927608108
$ws.Range("F11").Value = 1.019364922060666
$ws.Range("I11").Value = 1.025693496169847
$ws.Range("J11").Value = 1.015291165930345
$ws.Range("K11").Value = 1.015209756964894
$ws.Range("L11").Value = 1.014768835397501
$ws.Range("M11").Value = 1.022994810329474
$ws.Range("N11").Value = 1.01673299584149
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.008170768101115
$ws.Range("D12").Value = 1.011216436068369
$ws.Range("E12").Value = 1.010801719038985
$ws.Range("F12").Value = 1.018988759944971
$ws.Range("I12").Value = 1.025689670088849
$ws.Range("J12").Value = 1.01502590131766
$ws.Range("K12").Value = 1.014924161077123
$ws.Range("L12").Value = 1.014511088918054
$ws.Range("M12").Value = 1.022665933168186
$ws.Range("N12").Value = 1.016467354522629
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.008247117975416
$ws.Range("D13").Value = 1.011287902451848
$ws.Range("E13").Value = 1.010867192854159
$ws.Range("F13").Value = 1.019069471427794
$ws.Range("I13").Value = 1.025690512004175
$ws.Range("J13").Value = 1.015082812907864
$ws.Range("K13").Value = 1.014985432343158
$ws.Range("L13").Value = 1.014566386366101
$ws.Range("M13").Value = 1.022736505596643
$ws.Range("N13").Value = 1.01652434693382
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.008497248688079
$ws.Range("D14").Value = 1.011522045060138
$ws.Range("E14").Value = 1.011081701462394
$ws.Range("F14").Value = 1.019333838935127
$ws.Range("I14").Value = 1.025693189465947
$ws.Range("J14").Value = 1.015269244219557
$ws.Range("K14").Value = 1.015186154009716
$ws.Range("L14").Value = 1.014747534472642
$ws.Range("M14").Value = 1.022967637475574
$ws.Range("N14").Value = 1.016711042999359
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.008651335586056
$ws.Range("D15").Value = 1.011666291095762
$ws.Range("E15").Value = 1.011213851165168
$ws.Range("F15").Value = 1.019496655995158
$ws.Range("I15").Value = 1.025694777029423
$ws.Range("J15").Value = 1.0153840773264
$ws.Range("K15").Value = 1.015309796131759
$ws.Range("L15").Value = 1.01485911674117
$ws.Range("M15").Value = 1.023109966030379
$ws.Range("N15").Value = 1.016826039182385
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.009547783234548
$ws.Range("D16").Value = 1.012505611322986
$ws.Range("E16").Value = 1.011982777433198
$ws.Range("F16").Value = 1.020443281588223
$ws.Range("I16").Value = 1.025703068691744
$ws.Range("J16").Value = 1.016051950112728
$ws.Range("K16").Value = 1.016029006397017
$ws.Range("L16").Value = 1.015508133902475
$ws.Range("M16").Value = 1.023937171844764
$ws.Range("N16").Value = 1.017494860424676
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.010109751083007
$ws.Range("D17").Value = 1.013031876180644
$ws.Range("E17").Value = 1.012464895838376
$ws.Range("F17").Value = 1.021036164875737
$ws.Range("I17").Value = 1.025707433980311
$ws.Range("J17").Value = 1.016470448266561
$ws.Range("K17").Value = 1.016479764860676
$ws.Range("L17").Value = 1.015914861943557
$ws.Range("M17").Value = 1.024454997566233
$ws.Range("N17").Value = 1.01791395289389
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.010437413762242
$ws.Range("D18").Value = 1.013338760856178
$ws.Range("E18").Value = 1.012746033903416
$ws.Range("F18").Value = 1.021381658139332
$ws.Range("I18").Value = 1.025709678817867
$ws.Range("J18").Value = 1.016714394051073
$ws.Range("K18").Value = 1.016742548413554
$ws.Range("L18").Value = 1.016151963029983
$ws.Range("M18").Value = 1.024756657086487
$ws.Range("N18").Value = 1.018158245109397
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.010549117528
$ws.Range("D19").Value = 1.013443387779851
$ws.Range("E19").Value = 1.012841882427299
$ws.Range("F19").Value = 1.02149940742302
$ws.Range("I19").Value = 1.025710393120453
$ws.Range("J19").Value = 1.016797546774868
$ws.Range("K19").Value = 1.016832127879148
$ws.Range("L19").Value = 1.016232785373023
$ws.Range("M19").Value = 1.024859451015529
$ws.Range("N19").Value = 1.018241515919597
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.010049470098377
$ws.Range("D20").Value = 1.012975420922873
$ws.Range("E20").Value = 1.012413176733018
$ws.Range("F20").Value = 1.02097258783495
$ws.Range("I20").Value = 1.025706996794433
$ws.Range("J20").Value = 1.016425563665198
$ws.Range("K20").Value = 1.016431416876245
$ws.Range("L20").Value = 1.01587123802944
$ws.Range("M20").Value = 1.024399479062857
$ws.Range("N20").Value = 1.017869004551244
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.008423597043329
$ws.Range("D21").Value = 1.011453099441086
$ws.Range("E21").Value = 1.011018537394112
$ws.Range("F21").Value = 1.019256003574281
$ws.Range("I21").Value = 1.025692413958264
$ws.Range("J21").Value = 1.015214351808064
$ws.Range("K21").Value = 1.015127052558808
$ws.Range("L21").Value = 1.014694196933009
$ws.Range("M21").Value = 1.022899591455533
$ws.Range("N21").Value = 1.016656072634345
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.007400178378734
$ws.Range("D22").Value = 1.010495220460842
$ws.Range("E22").Value = 1.01014097247986
$ws.Range("F22").Value = 1.018173734693691
$ws.Range("I22").Value = 1.025680532987075
$ws.Range("J22").Value = 1.014451361920087
$ws.Range("K22").Value = 1.014305680081676
$ws.Range("L22").Value = 1.013952879820668
$ws.Range("M22").Value = 1.021953089920399
$ws.Range("N22").Value = 1.015891999213212
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.00794282527755
$ws.Range("D23").Value = 1.011003082043762
$ws.Range("E23").Value = 1.010606254435377
$ws.Range("F23").Value = 1.018747751287618
$ws.Range("I23").Value = 1.025687088264402
$ws.Range("J23").Value = 1.014855976842653
$ws.Range("K23").Value = 1.014741227180607
$ws.Range("L23").Value = 1.014345987623651
$ws.Range("M23").Value = 1.022455178663856
$ws.Range("N23").Value = 1.016297188735377
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.010076708862929
$ws.Range("D24").Value = 1.013000930861113
$ws.Range("E24").Value = 1.012436546597748
$ws.Range("F24").Value = 1.021001316569301
$ws.Range("I24").Value = 1.025707195271241
$ws.Range("J24").Value = 1.01644584557022
$ws.Range("K24").Value = 1.01645326367051
$ws.Range("L24").Value = 1.015890950220484
$ws.Range("M24").Value = 1.024424566657676
$ws.Range("N24").Value = 1.017889315258898
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.012548138142456
$ws.Range("D25").Value = 1.01531633602014
$ws.Range("E25").Value = 1.014557636859458
$ws.Range("F25").Value = 1.023603709497395
$ws.Range("I25").Value = 1.025718704832526
$ws.Range("J25").Value = 1.018284655533117
$ws.Range("K25").Value = 1.018434662437227
$ws.Range("L25").Value = 1.017678456280314
$ws.Range("M25").Value = 1.026695068112736
$ws.Range("N25").Value = 1.019730736542857
